$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 25703000
$ws.Range("B11").Value = "Ангарский"
$ws.Range("C11").Value = 2020
$ws.Range("D11").Value = -194
$ws.Range("E11").Value = 236912
$ws.Range("F11").Value = [double]"0.22602907408658068"
$ws.Range("G11").Value = [double]"34942.071599999996"
$ws.Range("H11").Value = [double]"0.99471533734044704"
$ws.Range("I11").Value = [double]"3.4392517052745326E-2"
$ws.Range("J11").Value = [double]"49.184704711622878"
$ws.Range("K11").Value = 23.7
$ws.Range("L11").Value = [double]"1.7052745323157966E-3"
$ws.Range("M11").Value = [double]"3.8664145336665093E-3"
$ws.Range("N11").Value = [double]"1.9935672317147298E-3"
$ws.Range("O11").Value = [double]"3.7435841156209904E-2"
$ws.Range("P11").Value = [double]"0.29477189842641993"
$ws.Range("Q11").Value = [double]"7.9265926504356052"
$ws.Range("R11").Value = [double]"7.5977578172485983E-5"
$ws.Range("S11").Value = [double]"5.778094820017559E-2"
$ws.Range("T11").Value = [double]"286.70030183089074"

$ws.Range("A11:E11").Style = $ws.Range("A10:E10").Style
$ws.Range("F11").Style = $ws.Range("F10").Style
$ws.Range("G11").Style = $ws.Range("G10").Style
$ws.Range("H11:T11").Style = $ws.Range("H10:T10").Style

$ws.Range("F16").Select()
